$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037409208979866
$ws.Range("D2").Value = 1.046455098560166
$ws.Range("E2").Value = 1.046329947878131
$ws.Range("F2").Value = 1.057534494387406
$ws.Range("I2").Value = 1.039860442664732
$ws.Range("J2").Value = 1.042512237786595
$ws.Range("K2").Value = 1.049220486956129
$ws.Range("L2").Value = 1.049095686882283
$ws.Range("M2").Value = 1.060269200266895
$ws.Range("N2").Value = 1.018110133954396
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038191488801445
$ws.Range("D3").Value = 1.047077638404125
$ws.Range("E3").Value = 1.04702604286292
$ws.Range("F3").Value = 1.058312061925468
$ws.Range("I3").Value = 1.040023259400307
$ws.Range("J3").Value = 1.042939725194611
$ws.Range("K3").Value = 1.049655164304568
$ws.Range("L3").Value = 1.049603702907099
$ws.Range("M3").Value = 1.060860715220925
$ws.Range("N3").Value = 1.018253032722273
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038698447108917
$ws.Range("D4").Value = 1.047481150600414
$ws.Range("E4").Value = 1.04747754360072
$ws.Range("F4").Value = 1.058816409259165
$ws.Range("I4").Value = 1.040127798568932
$ws.Range("J4").Value = 1.043216415866325
$ws.Range("K4").Value = 1.04993642322223
$ws.Range("L4").Value = 1.04993282513034
$ws.Range("M4").Value = 1.061244010717403
$ws.Range("N4").Value = 1.018345493054696
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038911755016195
$ws.Range("D5").Value = 1.04765094967952
$ws.Range("E5").Value = 1.047667611185007
$ws.Range("F5").Value = 1.059028724136747
$ws.Range("I5").Value = 1.040171551152774
$ws.Range("J5").Value = 1.04333275386415
$ws.Range("K5").Value = 1.050054661470816
$ws.Range("L5").Value = 1.050071282741632
$ws.Range("M5").Value = 1.061405276989244
$ws.Range("N5").Value = 1.018384361751371
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038947581006062
$ws.Range("D6").Value = 1.047679469166354
$ws.Range("E6").Value = 1.047699539360404
$ws.Range("F6").Value = 1.059064389519534
$ws.Range("I6").Value = 1.040178885906715
$ws.Range("J6").Value = 1.043352288496982
$ws.Range("K6").Value = 1.050074513974332
$ws.Range("L6").Value = 1.050094535886817
$ws.Range("M6").Value = 1.061432361808701
$ws.Range("N6").Value = 1.018390887866738
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038701296624404
$ws.Range("D7").Value = 1.047483418827802
$ws.Range("E7").Value = 1.047480082287118
$ws.Range("F7").Value = 1.058819245095004
$ws.Range("I7").Value = 1.040128383962963
$ws.Range("J7").Value = 1.043217970314326
$ws.Range("K7").Value = 1.04993800314113
$ws.Range("L7").Value = 1.049934674837883
$ws.Range("M7").Value = 1.061246165061226
$ws.Range("N7").Value = 1.01834601242714
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037673423314392
$ws.Range("D8").Value = 1.046665345313646
$ws.Range("E8").Value = 1.046564971428122
$ws.Range("F8").Value = 1.057797025286743
$ws.Range("I8").Value = 1.039915635451246
$ws.Range("J8").Value = 1.042656692238718
$ws.Range("K8").Value = 1.049367388641212
$ws.Range("L8").Value = 1.049267289111674
$ws.Range("M8").Value = 1.060468991420633
$ws.Range("N8").Value = 1.018158427942638
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035868169412326
$ws.Range("D9").Value = 1.045229146914693
$ws.Range("E9").Value = 1.044960797257601
$ws.Range("F9").Value = 1.056005101844934
$ws.Range("I9").Value = 1.039534544325248
$ws.Range("J9").Value = 1.041668304524733
$ws.Range("K9").Value = 1.048361909062666
$ws.Range("L9").Value = 1.048094418196847
$ws.Range("M9").Value = 1.059103772295319
$ws.Range("N9").Value = 1.01782786583959
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034668810021951
$ws.Range("D10").Value = 1.044275398372027
$ws.Range("E10").Value = 1.043897097142074
$ws.Range("F10").Value = 1.05481691082723
$ws.Range("I10").Value = 1.039276360791081
$ws.Range("J10").Value = 1.041009902431747
$ws.Range("K10").Value = 1.04769168816388
$ws.Range("L10").Value = 1.047314715584601
$ws.Range("M10").Value = 1.058196599026083
$ws.Range("N10").Value = 1.017607511156948
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034150481390803
$ws.Range("D11").Value = 1.043863321602796
$ws.Range("E11").Value = 1.043437893054589
$ws.Range("F11").Value = 1.054303962966243
$ws.Range("I11").Value = 1.039163596446128
$ws.Range("J11").Value = 1.040724948504851
$ws.Range("K11").Value = 1.047401517221029
$ws.Range("L11").Value = 1.046977640355315
$ws.Range("M11").Value = 1.05780451113282
$ws.Range("N11").Value = 1.017512106188111
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033958103571638
$ws.Range("D12").Value = 1.043710395614958
$ws.Range("E12").Value = 1.043267534441028
$ws.Range("F12").Value = 1.054113665964449
$ws.Range("I12").Value = 1.039121565927428
$ws.Range("J12").Value = 1.040619126150871
$ws.Range("K12").Value = 1.047293742205929
$ws.Range("L12").Value = 1.04685251854531
$ws.Range("M12").Value = 1.057658982898849
$ws.Range("N12").Value = 1.017476670556374
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033999362333401
$ws.Range("D13").Value = 1.043743192495481
$ws.Range("E13").Value = 1.043304067394265
$ws.Range("F13").Value = 1.054154474660971
$ws.Range("I13").Value = 1.039130588157134
$ws.Range("J13").Value = 1.040641824376791
$ws.Range("K13").Value = 1.04731685995367
$ws.Range("L13").Value = 1.046879353818884
$ws.Range("M13").Value = 1.057690194161001
$ws.Range("N13").Value = 1.017484271520011
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034134576257529
$ws.Range("D14").Value = 1.043850677872569
$ws.Range("E14").Value = 1.043423806857457
$ws.Range("F14").Value = 1.054288228156166
$ws.Range("I14").Value = 1.039160125140072
$ws.Range("J14").Value = 1.040716200734583
$ws.Range("K14").Value = 1.047392608345093
$ws.Range("L14").Value = 1.046967296046087
$ws.Range("M14").Value = 1.057792479455176
$ws.Range("N14").Value = 1.017509177023552
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03421790628779
$ws.Range("D15").Value = 1.043916921465288
$ws.Range("E15").Value = 1.043497610214109
$ws.Range("F15").Value = 1.054370669264151
$ws.Range("I15").Value = 1.039178304683589
$ws.Range("J15").Value = 1.040762029450158
$ws.Range("K15").Value = 1.047439280453973
$ws.Range("L15").Value = 1.047021491193816
$ws.Range("M15").Value = 1.057855515521353
$ws.Range("N15").Value = 1.017524522409386
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034703231064728
$ws.Range("D16").Value = 1.044302765743635
$ws.Range("E16").Value = 1.043927602369637
$ws.Range("F16").Value = 1.054850986270786
$ws.Range("I16").Value = 1.039283824237085
$ws.Range("J16").Value = 1.041028816924349
$ws.Range("K16").Value = 1.047710946807381
$ws.Range("L16").Value = 1.047337097698572
$ws.Range("M16").Value = 1.058222636024826
$ws.Range("N16").Value = 1.01761384312606
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035007932115551
$ws.Range("D17").Value = 1.044545038784089
$ws.Range("E17").Value = 1.044197697424343
$ws.Range("F17").Value = 1.055152692156601
$ws.Range("I17").Value = 1.039349754945579
$ws.Range("J17").Value = 1.041196203741126
$ws.Range("K17").Value = 1.047881367478722
$ws.Range("L17").Value = 1.047535215357448
$ws.Range("M17").Value = 1.058453116209887
$ws.Range("N17").Value = 1.01766987470597
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035185755618635
$ws.Range("D18").Value = 1.044686439503683
$ws.Range("E18").Value = 1.044355372797227
$ws.Range("F18").Value = 1.055328821174423
$ws.Range("I18").Value = 1.039388117648593
$ws.Range("J18").Value = 1.041293850910411
$ws.Range("K18").Value = 1.04798077464691
$ws.Range("L18").Value = 1.047650826139112
$ws.Range("M18").Value = 1.058587621143604
$ws.Range("N18").Value = 1.017702557915318
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035246405151392
$ws.Range("D19").Value = 1.044734668181138
$ws.Range("E19").Value = 1.044409158603032
$ws.Range("F19").Value = 1.055388901834312
$ws.Range("I19").Value = 1.039401182444983
$ws.Range("J19").Value = 1.041327148279788
$ws.Range("K19").Value = 1.048014670499427
$ws.Range("L19").Value = 1.047690255215774
$ws.Range("M19").Value = 1.058633495622282
$ws.Range("N19").Value = 1.01771370218847
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034975230577027
$ws.Range("D20").Value = 1.044519036170782
$ws.Range("E20").Value = 1.044168704957483
$ws.Range("F20").Value = 1.05512030653384
$ws.Range("I20").Value = 1.039342690876594
$ws.Range("J20").Value = 1.041178243336968
$ws.Range("K20").Value = 1.04786308256261
$ws.Range("L20").Value = 1.047513953806543
$ws.Range("M20").Value = 1.058428380640086
$ws.Range("N20").Value = 1.017663862947153
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034094754900904
$ws.Range("D21").Value = 1.043819022291071
$ws.Range("E21").Value = 1.04338854074943
$ws.Range("F21").Value = 1.054248834588973
$ws.Range("I21").Value = 1.039151431226205
$ws.Range("J21").Value = 1.040694298126235
$ws.Range("K21").Value = 1.047370302115561
$ws.Range("L21").Value = 1.046941396949432
$ws.Range("M21").Value = 1.057762355915414
$ws.Range("N21").Value = 1.017501842913616
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03354204827666
$ws.Range("D22").Value = 1.043379693326286
$ws.Range("E22").Value = 1.042899237932414
$ws.Range("F22").Value = 1.053702264469662
$ws.Range("I22").Value = 1.039030341198777
$ws.Range("J22").Value = 1.040390151832056
$ws.Range("K22").Value = 1.047060515012233
$ws.Range("L22").Value = 1.046581888394071
$ws.Range("M22").Value = 1.057344240951163
$ws.Range("N22").Value = 1.017399986436181
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033834964116903
$ws.Range("D23").Value = 1.043612513637865
$ws.Range("E23").Value = 1.043158510506601
$ws.Range("F23").Value = 1.053991882057767
$ws.Range("I23").Value = 1.039094612435052
$ws.Range("J23").Value = 1.040551372813399
$ws.Range("K23").Value = 1.047224734386289
$ws.Range("L23").Value = 1.046772424573593
$ws.Range("M23").Value = 1.057565830198769
$ws.Range("N23").Value = 1.017453981179451
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03499000669665
$ws.Range("D24").Value = 1.044530785364074
$ws.Range("E24").Value = 1.044181804993099
$ws.Range("F24").Value = 1.055134939742427
$ws.Range("I24").Value = 1.0393458831145
$ws.Range("J24").Value = 1.041186358829694
$ws.Range("K24").Value = 1.047871344717987
$ws.Range("L24").Value = 1.047523560825435
$ws.Range("M24").Value = 1.058439557364114
$ws.Range("N24").Value = 1.017666579399363
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03633414903054
$ws.Range("D25").Value = 1.045599792165277
$ws.Range("E25").Value = 1.0453745103771
$ws.Range("F25").Value = 1.056467234319994
$ws.Range("I25").Value = 1.039633795592852
$ws.Range("J25").Value = 1.041923740252614
$ws.Range("K25").Value = 1.048621838658971
$ws.Range("L25").Value = 1.048397250615456
$ws.Range("M25").Value = 1.059456198127276
$ws.Range("N25").Value = 1.01791332270132
